# F16 input.xlsx update
#  - "Pmax"/"Pmax-Pin" block moved from Parameter sheet rows 25/26 to rows 23/24,
#    with the Data-sheet reference switched from B2 to E2 and a couple of values
#    changed (CO2 price row etc.)
#  - Data sheet header row relabeled + columns (CO2 price / interest rate /
#    renewable factor / total demand) reordered/shuffled, and its own
#    "Pmax" summary line moved from row 5 to row 4 (col B/C -> col D/E).
#  - "prices and emission factors" sheet: drop the stray "Nan" label in C3 and
#    a couple of now-unused columns.
#  - a batch of capacity values on the Parameter sheet changed.
#  - header text relabeled for 3 shared strings (energy carrier / prices / emission factor).
#  - new orange fill + a couple of style tweaks for the Pmax/Pmax-Pin cells.

$wb = $excel.ActiveWorkbook

$wsParams = $wb.Worksheets.Item("Parameter for Powerplants")
$wsPrices = $wb.Worksheets.Item("prices and emmision factors")
$wsData   = $wb.Worksheets.Item("Data")

# ---------------------------------------------------------------------------
# 1. "Parameter for Powerplants" sheet — installed capacity tweaks
# ---------------------------------------------------------------------------
$wsParams.Range("C2").Value = 2500
$wsParams.Range("C3").Value = 2500
$wsParams.Range("K3").Value = 0.1
$wsParams.Range("C5").Value = 0
$wsParams.Range("C6").Value = 0
$wsParams.Range("C9").Value = 0
$wsParams.Range("C12").Value = 0
$wsParams.Range("C15").Value = 0

# Replace the old "Pmax" (row25) / "Pmax-sum(P)" (row26) summary block with the
# new "Pmax" (row24) / "Pmax-Pin" (row23) block, referencing Data!E2 now.
$wsParams.Range("A22:K26").Clear()

$wsParams.Range("C24").Formula = "=+(0.000606221407201*Data!E2)"
$wsParams.Range("D24").Value = "Pmax"

$wsParams.Range("C23").Formula = "=C24-SUM(C2:C21)"
$wsParams.Range("D23").Value = "Pmax-Pin"

# Style: both new label/value cells get bold font + thin border + a solid
# fill (yellow for the Pmax row, new orange for the Pmax-Pin row), with no
# explicit alignment set. (NOTE: apply per-cell — a multi-area "C24,D24"
# Range only actually formats its first area.)
foreach ($addr in @("C24", "D24")) {
    $r = $wsParams.Range($addr)
    $r.Font.Bold = $true
    $r.Interior.Color = 65535
    $r.Borders.LineStyle = 1
    $r.Borders.Weight = 2
}

foreach ($addr in @("C23", "D23")) {
    $r = $wsParams.Range($addr)
    $r.Font.Bold = $true
    $r.Interior.Color = 49407
    $r.Borders.LineStyle = 1
    $r.Borders.Weight = 2
}

# ---------------------------------------------------------------------------
# 2. "prices and emmision factors" sheet
# ---------------------------------------------------------------------------
$wsPrices.Range("B1").Value = "energy carrier"
$wsPrices.Range("C1").Value = "prices(EUR/MWh)"
$wsPrices.Range("D1").Value = "emission factor"
$wsPrices.Range("C3").ClearContents()

# ---------------------------------------------------------------------------
# 3. "Data" sheet — relabel header row + reorder the stored values, then
#    rebuild the "Pmax" helper line on row 4 (was row 5) using columns D/E.
# ---------------------------------------------------------------------------
$wsData.Range("B1").Value = "CO2 Price"
$wsData.Range("C1").Value = "Interes Rate [0-1]"
$wsData.Range("D1").Value = "Total Renewable Factor [0-1]"
$wsData.Range("E1").Value = "Total Demand[ MWh]"

$wsData.Range("B2").Value = 25
$wsData.Range("C2").Value = 0.5
$wsData.Range("D2").Value = 0
$wsData.Range("E2").Value = 6000000

$wsData.Range("A3:E5").Clear()

$wsData.Range("D4").Value = "Pmax"
$wsData.Range("E4").Formula = "=0.000606221407201*E2"

foreach ($addr in @("D4", "E4")) {
    $r = $wsData.Range($addr)
    $r.Font.Bold = $true
    $r.Interior.Color = 65535
    $r.Borders.LineStyle = 1
    $r.Borders.Weight = 2
}

# D4 (the "Pmax" label) is right/center aligned; E4 (the formula result) keeps
# the default (no explicit) alignment.
$wsData.Range("D4").HorizontalAlignment = -4152
$wsData.Range("D4").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 4. Selections / active sheet — "Data" ends up the active/visible tab.
# ---------------------------------------------------------------------------
$wsParams.Activate()
$wsParams.Range("C23").Select()

$wsPrices.Activate()
$wsPrices.Range("C25").Select()

$wsData.Activate()
$wsData.Range("E15").Select()
